$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.013704583048820496
$ws.Range("C2").Value = 0.006486182566732168
$ws.Range("D2").Value = 0.00432577496394515
$ws.Range("E2").Value = 0.004177039489150047
$ws.Range("F2").Value = 0.00015976434224285185
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.12727367877960205
$ws.Range("K2").Value = 1.415645718574524
